$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.619.09"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.316.99"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "576.36"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "173.99"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "3.313.83"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "45.68"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "700.14"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "3.862.86"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "8.37"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "67.642.49"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "3.321.71"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").Value = "17.35"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "0.887"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "5.36"
$ws.Range("E23").Value = "  +3.65%  "
$ws.Range("D24").Value = "16.87"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "97.66"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("D28").Value = "9.34"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "33.03"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").Value = "8.45"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "7.09"
$ws.Range("E31").Value = "  +5.32%  "
$ws.Range("D32").Value = "566.73"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("D33").Value = "10.92"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "57.47"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "3.712.69"
$ws.Range("E37").Value = "  -4.59%  "
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").Value = "34.17"
$ws.Range("E39").Value = "  +5.70%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "3.16"
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").Value = "3.31"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").Value = "0.332"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "0.0₃0665"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  +7.39%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0403"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "1.31"
$ws.Range("E50").Value = "  -5.30%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "128.02"
$ws.Range("E51").Value = "  -1.91%  "
